$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(181, "Sunday, Jan 15", "4:10 PM", "FR4044", "Liverpool", "(LPL)", "Ryanair ", "B738", "(SP-RKP)", "4:15 PM", "0 hours, 5 minutes"),
    @(182, "Sunday, Jan 15", "5:50 PM", "FR3284", "Riga",      "(RIX)", "Buzz ",    "B38M", "(SP-RZE)", "5:57 PM", "0 hours, 7 minutes"),
    @(183, "Sunday, Jan 15", "6:05 PM", "FR1056", "Brussels",  "(CRL)", "Ryanair ", "B738", "(SP-RSP)", "6:05 PM", "0 hours, 0 minutes")
)

$r = 182
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Range("K$r").ClearFormats()
    $ws.Cells.Item($r, 12).Value = $row[10]
    $ws.Range("M$r").ClearFormats()
    $r = $r + 1
}
